# Add "Goldsmith 2003 JCR" row to the Empirical Strategies Summary sheet.
# The table is sorted alphabetically by Authors (column A); "Goldsmith"
# sits between "Gates and Tersawa" (row 8) and "Hansen et al" (row 9),
# so insert a new row at 9 and shift everything else down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "Goldsmith"
$ws.Range("B9").Value = 2003
$ws.Range("C9").Value = "JCR"
$ws.Range("E9").Value = "OLS: PCSE"
$ws.Range("F9").Value = "Expendityres / GDP"
$ws.Range("G9").Value = "Dummy: Defense Pact"
$ws.Range("H9").Value = "None"
$ws.Range("I9").Value = "1886-1989"
$ws.Range("J9").Value = "All States"
$ws.Range("K9").Value = "GDP, regime type, lagged DV, war, rivalries, regional context, major power, systemic variables"
$ws.Range("L9").Value = "General"

# Restore the view to the top-left / default zoom with the new selection,
# matching the saved workbook state after the edit.
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("N9").Select()
